# Scheduled runner update: refresh market price / profit figures
# (currentAveragePrice*, LevePrice*, LeveProfit* columns H:N) across the
# various crafting-job sheets, as produced by the automated data sync.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1425
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 1350
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 1350
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -1700

$ws.Range("H64").Value = 3000
$ws.Range("I64").Value = 3000
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 3000
$ws.Range("L64").Value = 3000
$ws.Range("M64").Value = -2752
$ws.Range("N64").Value = -3496

$ws.Range("H67").Value = 3000
$ws.Range("I67").Value = 3000
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 3000
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = -2142
$ws.Range("N67").Value = -4716

$ws.Range("H76").Value = 3250.375
$ws.Range("I76").Value = 3267.1667
$ws.Range("J76").Value = 3200
$ws.Range("K76").Value = 3267.1667
$ws.Range("L76").Value = 3200
$ws.Range("M76").Value = -2952.1667
$ws.Range("N76").Value = -3830

$ws.Range("H79").Value = 3250.375
$ws.Range("I79").Value = 3267.1667
$ws.Range("J79").Value = 3200
$ws.Range("K79").Value = 3267.1667
$ws.Range("L79").Value = 3200
$ws.Range("M79").Value = -2175.1667
$ws.Range("N79").Value = -5384

$ws.Range("H137").Value = 2033945.6
$ws.Range("I137").Value = 4387067.5
$ws.Range("J137").Value = 1703.8636
$ws.Range("K137").Value = 13161202.5
$ws.Range("L137").Value = 5111.5908
$ws.Range("M137").Value = -13158652.5
$ws.Range("N137").Value = -10211.5908

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5962175
$ws.Range("I32").Value = 6674915.5
$ws.Range("J32").Value = 22666.666
$ws.Range("K32").Value = 6674915.5
$ws.Range("L32").Value = 22666.666
$ws.Range("M32").Value = -6674628.5
$ws.Range("N32").Value = -23240.666

$ws.Range("H63").Value = 4079.8096
$ws.Range("I63").Value = 3338.9
$ws.Range("J63").Value = 4753.364
$ws.Range("K63").Value = 3338.9
$ws.Range("L63").Value = 4753.364
$ws.Range("M63").Value = -2652.9
$ws.Range("N63").Value = -6125.364

$ws.Range("H66").Value = 4079.8096
$ws.Range("I66").Value = 3338.9
$ws.Range("J66").Value = 4753.364
$ws.Range("K66").Value = 16694.5
$ws.Range("L66").Value = 23766.82
$ws.Range("M66").Value = -13262.5
$ws.Range("N66").Value = -30630.82

$ws.Range("H74").Value = 23812782
$ws.Range("I74").Value = 2743.4167
$ws.Range("J74").Value = 55559500
$ws.Range("K74").Value = 2743.4167
$ws.Range("L74").Value = 55559500
$ws.Range("M74").Value = -1869.4167
$ws.Range("N74").Value = -55561248

$ws.Range("H77").Value = 23812782
$ws.Range("I77").Value = 2743.4167
$ws.Range("J77").Value = 55559500
$ws.Range("K77").Value = 13717.0835
$ws.Range("L77").Value = 277797500
$ws.Range("M77").Value = -9349.083500000001
$ws.Range("N77").Value = -277806236

$ws.Range("H132").Value = 988355.4
$ws.Range("I132").Value = 1578.459
$ws.Range("J132").Value = 4529143
$ws.Range("K132").Value = 4735.377
$ws.Range("L132").Value = 13587429
$ws.Range("M132").Value = -2205.377
$ws.Range("N132").Value = -13592489

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 31251242
$ws.Range("I105").Value = 31251242
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 31251242
$ws.Range("L105").Value = 0
$ws.Range("N105").Value = -31249495
$ws.Range("M105").ClearContents()

$ws.Range("H134").Value = 2370.5264
$ws.Range("I134").Value = 2280
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 6840
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -4305
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5975.4
$ws.Range("I31").Value = 1777.125
$ws.Range("J31").Value = 7219.3335
$ws.Range("K31").Value = 1777.125
$ws.Range("L31").Value = 7219.3335
$ws.Range("M31").Value = -1482.125
$ws.Range("N31").Value = -7809.3335

$ws.Range("H34").Value = 5975.4
$ws.Range("I34").Value = 1777.125
$ws.Range("J34").Value = 7219.3335
$ws.Range("K34").Value = 1777.125
$ws.Range("L34").Value = 7219.3335
$ws.Range("M34").Value = -1575.125
$ws.Range("N34").Value = -7623.3335

$ws.Range("H62").Value = 4139.6895
$ws.Range("I62").Value = 3987.3125
$ws.Range("J62").Value = 4327.231
$ws.Range("K62").Value = 3987.3125
$ws.Range("L62").Value = 4327.231
$ws.Range("M62").Value = -3363.3125
$ws.Range("N62").Value = -5575.231

$ws.Range("H65").Value = 4139.6895
$ws.Range("I65").Value = 3987.3125
$ws.Range("J65").Value = 4327.231
$ws.Range("K65").Value = 19936.5625
$ws.Range("L65").Value = 21636.155
$ws.Range("M65").Value = -16816.5625
$ws.Range("N65").Value = -27876.155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6088293.5
$ws.Range("I4").Value = 11667137
$ws.Range("J4").Value = 2282
$ws.Range("K4").Value = 35001411
$ws.Range("L4").Value = 6846
$ws.Range("M4").Value = -35001299
$ws.Range("N4").Value = -7070

$ws.Range("H15").Value = 524.8261
$ws.Range("I15").Value = 690.3333
$ws.Range("J15").Value = 500
$ws.Range("K15").Value = 2070.9999
$ws.Range("L15").Value = 1500
$ws.Range("M15").Value = -1930.9999
$ws.Range("N15").Value = -1780

$ws.Range("H16").Value = 300
$ws.Range("I16").Value = 300
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 900
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -727

$ws.Range("H20").Value = 976.2083
$ws.Range("I20").Value = 810
$ws.Range("J20").Value = 999.9524
$ws.Range("K20").Value = 2430
$ws.Range("L20").Value = 2999.8572
$ws.Range("M20").Value = -2203
$ws.Range("N20").Value = -3453.8572

$ws.Range("H21").Value = 2240
$ws.Range("I21").Value = 200
$ws.Range("J21").Value = 3600
$ws.Range("K21").Value = 600
$ws.Range("L21").Value = 10800
$ws.Range("M21").Value = -427
$ws.Range("N21").Value = -11146

$ws.Range("H113").Value = 556.4375
$ws.Range("I113").Value = 532.1
$ws.Range("J113").Value = 597
$ws.Range("K113").Value = 1596.3
$ws.Range("L113").Value = 1791
$ws.Range("M113").Value = 573.6999999999998
$ws.Range("N113").Value = -6131

$ws.Range("H121").Value = 838.9231
$ws.Range("I121").Value = 294.7
$ws.Range("J121").Value = 2653
$ws.Range("K121").Value = 884.0999999999999
$ws.Range("L121").Value = 7959
$ws.Range("M121").Value = 425.9000000000001
$ws.Range("N121").Value = -10579

$ws.Range("H131").Value = 4271.1943
$ws.Range("I131").Value = 643.75
$ws.Range("J131").Value = 5307.607
$ws.Range("K131").Value = 1931.25
$ws.Range("L131").Value = 15922.821
$ws.Range("M131").Value = 3108.75
$ws.Range("N131").Value = -26002.821

$ws.Range("H134").Value = 5945.5474
$ws.Range("I134").Value = 1983.091
$ws.Range("J134").Value = 8757.612999999999
$ws.Range("K134").Value = 5949.272999999999
$ws.Range("L134").Value = 26272.839
$ws.Range("M134").Value = -879.2729999999992
$ws.Range("N134").Value = -36412.839

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1046633.3
$ws.Range("I80").Value = 1502216.6
$ws.Range("J80").Value = 135466.67
$ws.Range("K80").Value = 1502216.6
$ws.Range("L80").Value = 135466.67
$ws.Range("M80").Value = -1501218.6
$ws.Range("N80").Value = -137462.67

$ws.Range("H83").Value = 1046633.3
$ws.Range("I83").Value = 1502216.6
$ws.Range("J83").Value = 135466.67
$ws.Range("K83").Value = 7511083
$ws.Range("L83").Value = 677333.3500000001
$ws.Range("M83").Value = -7506091
$ws.Range("N83").Value = -687317.3500000001

$ws.Range("H132").Value = 24394468
$ws.Range("I132").Value = 29415946
$ws.Range("J132").Value = 4430
$ws.Range("K132").Value = 88247838
$ws.Range("L132").Value = 13290
$ws.Range("M132").Value = -88245308
$ws.Range("N132").Value = -18350

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5632.5
$ws.Range("I61").Value = 2000
$ws.Range("J61").Value = 6359
$ws.Range("K61").Value = 2000
$ws.Range("L61").Value = 6359
$ws.Range("M61").Value = -1798
$ws.Range("N61").Value = -6763

$ws.Range("H68").Value = 2250
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 2250
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 2250
$ws.Range("N68").Value = -3748
$ws.Range("M68").ClearContents()

$ws.Range("H71").Value = 2250
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 2250
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 11250
$ws.Range("N71").Value = -18738
$ws.Range("M71").ClearContents()

$ws.Range("H113").Value = 5632.5
$ws.Range("I113").Value = 2000
$ws.Range("J113").Value = 6359
$ws.Range("K113").Value = 2000
$ws.Range("L113").Value = 6359
$ws.Range("M113").Value = 170
$ws.Range("N113").Value = -10699

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6946622.5
$ws.Range("I132").Value = 2014.6061
$ws.Range("J132").Value = 32410186
$ws.Range("K132").Value = 6043.8183
$ws.Range("L132").Value = 97230558
$ws.Range("M132").Value = -3513.8183
$ws.Range("N132").Value = -97235618
